$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 20. This pushes the existing row 20 (Especial, 10kg tray)
# down to row 21, and the existing row 21 (Segunda, 18kg tray) down to row 22,
# leaving a fresh, empty row 20 (with formatting carried over) ready for the
# new weekly price entry.
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with this week's price record.
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = 45034
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100101
$ws.Cells.Item(20, 8).Value = "Berries"
$ws.Cells.Item(20, 9).Value = 100101007
$ws.Cells.Item(20, 10).Value = "Kiwi"
$ws.Cells.Item(20, 11).Value = "Hayward"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 250
$ws.Cells.Item(20, 14).Value = 25000
$ws.Cells.Item(20, 15).Value = 26000
$ws.Cells.Item(20, 16).Value = 25600
$ws.Cells.Item(20, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 1422
$ws.Cells.Item(20, 20).Value = 18
